# Update "想去人数" (interested-people count) figures scraped at a later time.
# Applies to both the "展览" sheet and the "全部类型" sheet (F column values).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# row -> new value for column F
$updates = @{
    5  = 13143
    12 = 13770
    13 = 14370
    21 = 35
    25 = 5425
    28 = 317
    30 = 49
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
